$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1, 2, 3, 4, 5, 6)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws.Range("A7").Select() | Out-Null
